$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.399.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7020"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.22%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07918"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3132"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07835"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.894.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "93.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.177"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7008"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.524"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008388"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.430.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.119.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.20%  "

$ws.Range("E21").Value = "  -1.28%  "

$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.648"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.62%  "

$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1559"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.005"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.79%  "

$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.314"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.260"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("E32").Value = "  +2.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05268"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.897"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.18%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7502"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.92%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.176"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.712"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.84%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.269.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.770"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8902"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.41%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.99%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.015"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.38%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000127"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.022.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.602"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.799"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.5183"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "

$ws.Range("E51").Value = "  -1.36%  "
